$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("gunslinger")
$ws1.Activate()
